$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column D (Price) and Column E (Volume(1h)) updates per row.
# Each entry: row, newD (or $null to skip), newE (or $null to skip).
# D values that look like plain decimal numbers are prefixed with a
# leading apostrophe so Excel keeps them as text (matching the
# original inline-string cell type) instead of silently converting
# them to floating point numbers.
$updates = @(
    @{Row=2;  D="41.525.58";   E=$null},
    @{Row=3;  D="2.474.16";    E="  +0.46%  "},
    @{Row=4;  D=$null;         E="  -0.04%  "},
    @{Row=5;  D="'312.74";     E="  +0.10%  "},
    @{Row=6;  D="'91.71";      E="  -2.52%  "},
    @{Row=7;  D="'0.547";      E="  +0.13%  "},
    @{Row=8;  D=$null;         E="  -0.14%  "},
    @{Row=9;  D=$null;         E="  +2.48%  "},
    @{Row=10; D="'32.52";      E="  -2.71%  "},
    @{Row=11; D="'0.0787";     E="  +0.85%  "},
    @{Row=12; D=$null;         E="  +1.14%  "},
    @{Row=13; D=$null;         E="  +0.54%  "},
    @{Row=14; D="'6.87";       E="  -1.47%  "},
    @{Row=15; D="'16.26";      E="  +8.85%  "},
    @{Row=16; D="2.451.93";    E="  -0.65%  "},
    @{Row=17; D=$null;         E="  -2.02%  "},
    @{Row=18; D="41.503.54";   E="  +0.66%  "},
    @{Row=19; D="'6.50";       E="  +3.00%  "},
    @{Row=20; D="0.0₃0940";    E="  +2.11%  "},
    @{Row=21; D="'71.89";      E="  +5.14%  "},
    @{Row=22; D="'11.05";      E="  -1.41%  "},
    @{Row=23; D="'236.29";     E="  -0.62%  "},
    @{Row=24; D=$null;         E="  -1.21%  "},
    @{Row=25; D=$null;         E="  -0.07%  "},
    @{Row=26; D="'1.90";       E="  +0.29%  "},
    @{Row=27; D="'24.90";      E="  +4.09%  "},
    @{Row=28; D=$null;         E="  -0.76%  "},
    @{Row=29; D="'9.68";       E="  +0.21%  "},
    @{Row=30; D="'35.69";      E="  -1.74%  "},
    @{Row=31; D="'157.06";     E="  +3.65%  "},
    @{Row=32; D="'5.45";       E="  -0.59%  "},
    @{Row=33; D=$null;         E="  -0.56%  "},
    @{Row=34; D=$null;         E="  +1.56%  "},
    @{Row=35; D="'17.35";      E="  +0.60%  "},
    @{Row=36; D=$null;         E="  -8.42%  "},
    @{Row=37; D=$null;         E="  -5.62%  "},
    @{Row=38; D=$null;         E="  +2.82%  "},
    @{Row=39; D="'1.82";       E="  -2.73%  "},
    @{Row=40; D=$null;         E="  -0.17%  "},
    @{Row=41; D=$null;         E="  -4.70%  "},
    @{Row=42; D=$null;         E="  -0.20%  "},
    @{Row=43; D="1.960.77";    E="  -1.03%  "},
    @{Row=44; D=$null;         E="  -0.25%  "},
    @{Row=45; D="'18.78";      E="  -4.15%  "},
    @{Row=46; D=$null;         E="  -2.40%  "},
    @{Row=47; D="'8.95";       E="  +3.12%  "},
    @{Row=48; D="2.716.45";    E="  +0.44%  "},
    @{Row=49; D="'97.59";      E="  +1.14%  "},
    @{Row=50; D="'67.38";      E="  -3.07%  "},
    @{Row=51; D="'72.08";      E=$null}
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}

$wb.Save()
